# Auto-generated edit script applying scheduled market-data refresh
# to the Leve profit sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 4138.533
$ws.Cells.Item(62, 9).Value = 2618.6
$ws.Cells.Item(62, 10).Value = 7178.4
$ws.Cells.Item(62, 11).Value = 2618.6
$ws.Cells.Item(62, 12).Value = 7178.4
$ws.Cells.Item(62, 13).Value = -1994.6
$ws.Cells.Item(62, 14).Value = -8426.4

$ws.Cells.Item(65, 8).Value = 4138.533
$ws.Cells.Item(65, 9).Value = 2618.6
$ws.Cells.Item(65, 10).Value = 7178.4
$ws.Cells.Item(65, 11).Value = 13093
$ws.Cells.Item(65, 12).Value = 35892
$ws.Cells.Item(65, 13).Value = -9973
$ws.Cells.Item(65, 14).Value = -42132

$ws.Cells.Item(129, 8).Value = 170389.72
$ws.Cells.Item(129, 10).Value = 209381.77
$ws.Cells.Item(129, 12).Value = 628145.3099999999
$ws.Cells.Item(129, 14).Value = -638145.3099999999

$ws.Cells.Item(132, 8).Value = 2558.2703
$ws.Cells.Item(132, 9).Value = 2558.2703
$ws.Cells.Item(132, 11).Value = 7674.8109
$ws.Cells.Item(132, 13).Value = -5144.8109

$ws.Cells.Item(137, 8).Value = 1319.5454
$ws.Cells.Item(137, 9).Value = 1067.2916
$ws.Cells.Item(137, 11).Value = 3201.8748
$ws.Cells.Item(137, 13).Value = -651.8748000000001

$ws.Cells.Item(138, 8).Value = 2180.291
$ws.Cells.Item(138, 9).Value = 1426.0526
$ws.Cells.Item(138, 10).Value = 2578.361
$ws.Cells.Item(138, 11).Value = 4278.1578
$ws.Cells.Item(138, 12).Value = 7735.083
$ws.Cells.Item(138, 13).Value = 861.8422
$ws.Cells.Item(138, 14).Value = -18015.083

$ws.Cells.Item(141, 8).Value = 2906.7273
$ws.Cells.Item(141, 9).Value = 2639.1428
$ws.Cells.Item(141, 11).Value = 7917.428400000001
$ws.Cells.Item(141, 13).Value = -2737.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 764.4400000000001
$ws.Cells.Item(2, 9).Value = 733.7917
$ws.Cells.Item(2, 10).Value = 1500
$ws.Cells.Item(2, 11).Value = 733.7917
$ws.Cells.Item(2, 12).Value = 1500
$ws.Cells.Item(2, 13).Value = -620.7917
$ws.Cells.Item(2, 14).Value = -1726

$ws.Cells.Item(32, 8).Value = 8588.125
$ws.Cells.Item(32, 9).Value = 7227.521
$ws.Cells.Item(32, 10).Value = 16751.75
$ws.Cells.Item(32, 11).Value = 7227.521
$ws.Cells.Item(32, 12).Value = 16751.75
$ws.Cells.Item(32, 13).Value = -6940.521
$ws.Cells.Item(32, 14).Value = -17325.75

$ws.Cells.Item(61, 8).Value = 1566.7567
$ws.Cells.Item(61, 9).Value = 1192.5807
$ws.Cells.Item(61, 11).Value = 1192.5807
$ws.Cells.Item(61, 13).Value = -980.5807

$ws.Cells.Item(63, 8).Value = 10418417
$ws.Cells.Item(63, 9).Value = 2625
$ws.Cells.Item(63, 11).Value = 2625
$ws.Cells.Item(63, 13).Value = -1939

$ws.Cells.Item(66, 8).Value = 10418417
$ws.Cells.Item(66, 9).Value = 2625
$ws.Cells.Item(66, 11).Value = 13125
$ws.Cells.Item(66, 13).Value = -9693

$ws.Cells.Item(74, 8).Value = 32259842
$ws.Cells.Item(74, 9).Value = 47619652
$ws.Cells.Item(74, 11).Value = 47619652
$ws.Cells.Item(74, 13).Value = -47618778

$ws.Cells.Item(77, 8).Value = 32259842
$ws.Cells.Item(77, 9).Value = 47619652
$ws.Cells.Item(77, 11).Value = 238098260
$ws.Cells.Item(77, 13).Value = -238093892

$ws.Cells.Item(116, 8).Value = 764.4400000000001
$ws.Cells.Item(116, 9).Value = 733.7917
$ws.Cells.Item(116, 10).Value = 1500
$ws.Cells.Item(116, 11).Value = 733.7917
$ws.Cells.Item(116, 12).Value = 1500
$ws.Cells.Item(116, 13).Value = 1560.2083
$ws.Cells.Item(116, 14).Value = -6088

$ws.Cells.Item(136, 8).Value = 1566.7567
$ws.Cells.Item(136, 9).Value = 1192.5807
$ws.Cells.Item(136, 11).Value = 3577.7421
$ws.Cells.Item(136, 13).Value = -1027.7421

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 764.4400000000001
$ws.Cells.Item(3, 9).Value = 733.7917
$ws.Cells.Item(3, 10).Value = 1500
$ws.Cells.Item(3, 11).Value = 733.7917
$ws.Cells.Item(3, 12).Value = 1500
$ws.Cells.Item(3, 13).Value = -619.7917
$ws.Cells.Item(3, 14).Value = -1728

$ws.Cells.Item(86, 8).Value = 1965.7858
$ws.Cells.Item(86, 9).Value = 1700.2941
$ws.Cells.Item(86, 10).Value = 2376.0908
$ws.Cells.Item(86, 11).Value = 1700.2941
$ws.Cells.Item(86, 12).Value = 2376.0908
$ws.Cells.Item(86, 13).Value = -577.2941000000001
$ws.Cells.Item(86, 14).Value = -4622.0908

$ws.Cells.Item(89, 8).Value = 1965.7858
$ws.Cells.Item(89, 9).Value = 1700.2941
$ws.Cells.Item(89, 10).Value = 2376.0908
$ws.Cells.Item(89, 11).Value = 8501.470499999999
$ws.Cells.Item(89, 12).Value = 11880.454
$ws.Cells.Item(89, 13).Value = -2885.470499999999
$ws.Cells.Item(89, 14).Value = -23112.454

$ws.Cells.Item(94, 8).Value = 787.4
$ws.Cells.Item(94, 9).Value = 556.2632
$ws.Cells.Item(94, 10).Value = 1061.875
$ws.Cells.Item(94, 11).Value = 556.2632
$ws.Cells.Item(94, 12).Value = 1061.875
$ws.Cells.Item(94, 13).Value = -105.2632
$ws.Cells.Item(94, 14).Value = -1963.875

$ws.Cells.Item(134, 8).Value = 3940.3845
$ws.Cells.Item(134, 9).Value = 4232.7
$ws.Cells.Item(134, 10).Value = 2966
$ws.Cells.Item(134, 11).Value = 12698.1
$ws.Cells.Item(134, 12).Value = 8898
$ws.Cells.Item(134, 13).Value = -10163.1
$ws.Cells.Item(134, 14).Value = -13968

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3104.7708
$ws.Cells.Item(31, 9).Value = 1607.2413
$ws.Cells.Item(31, 10).Value = 5390.4736
$ws.Cells.Item(31, 11).Value = 1607.2413
$ws.Cells.Item(31, 12).Value = 5390.4736
$ws.Cells.Item(31, 13).Value = -1312.2413
$ws.Cells.Item(31, 14).Value = -5980.4736

$ws.Cells.Item(34, 8).Value = 3104.7708
$ws.Cells.Item(34, 9).Value = 1607.2413
$ws.Cells.Item(34, 10).Value = 5390.4736
$ws.Cells.Item(34, 11).Value = 1607.2413
$ws.Cells.Item(34, 12).Value = 5390.4736
$ws.Cells.Item(34, 13).Value = -1405.2413
$ws.Cells.Item(34, 14).Value = -5794.4736

$ws.Cells.Item(58, 8).Value = 33174.438
$ws.Cells.Item(58, 10).Value = 85154.664
$ws.Cells.Item(58, 12).Value = 85154.664
$ws.Cells.Item(58, 14).Value = -85560.664

$ws.Cells.Item(132, 8).Value = 2018.7174
$ws.Cells.Item(132, 9).Value = 1584.1666
$ws.Cells.Item(132, 10).Value = 3583.1
$ws.Cells.Item(132, 11).Value = 4752.4998
$ws.Cells.Item(132, 12).Value = 10749.3
$ws.Cells.Item(132, 13).Value = -2222.4998
$ws.Cells.Item(132, 14).Value = -15809.3

$ws.Cells.Item(134, 8).Value = 825.1111
$ws.Cells.Item(134, 9).Value = 676.8
$ws.Cells.Item(134, 10).Value = 1566.6666
$ws.Cells.Item(134, 11).Value = 2030.4
$ws.Cells.Item(134, 12).Value = 4699.9998
$ws.Cells.Item(134, 13).Value = 504.6000000000001
$ws.Cells.Item(134, 14).Value = -9769.9998

$ws.Cells.Item(136, 8).Value = 33174.438
$ws.Cells.Item(136, 10).Value = 85154.664
$ws.Cells.Item(136, 12).Value = 255463.992
$ws.Cells.Item(136, 14).Value = -260563.992

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 699.16
$ws.Cells.Item(131, 10).Value = 723
$ws.Cells.Item(131, 12).Value = 2169
$ws.Cells.Item(131, 14).Value = -12249

$ws.Cells.Item(136, 8).Value = 2962.125
$ws.Cells.Item(136, 10).Value = 4991.75
$ws.Cells.Item(136, 12).Value = 14975.25
$ws.Cells.Item(136, 14).Value = -25175.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 9187.77
$ws.Cells.Item(113, 9).Value = 10860.1
$ws.Cells.Item(113, 11).Value = 10860.1
$ws.Cells.Item(113, 13).Value = -8690.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 14).ClearContents()

$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 12).Value = 0
$ws.Cells.Item(67, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 844.381
$ws.Cells.Item(113, 9).Value = 919.0526
$ws.Cells.Item(113, 11).Value = 2757.1578
$ws.Cells.Item(113, 13).Value = -587.1578
